$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.78
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 4.1
$ws.Range("U2").Value = 2.4
$ws.Range("W2").Value = 2.16
$ws.Range("K3").Value = 5.3
$ws.Range("N3").Value = 5.5
$ws.Range("Q3").Value = 1.64
$ws.Range("S3").Value = 2.6
$ws.Range("AA3").Value = 260
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 5.7
$ws.Range("F4").Value = 2.04
$ws.Range("P4").Value = 2.44
$ws.Range("S4").Value = 2.66
$ws.Range("W4").Value = 1.94
$ws.Range("AI4").Value = 40
$ws.Range("AO4").Value = 28
$ws.Range("F5").Value = 1.27
$ws.Range("I5").Value = 7.2
$ws.Range("J5").Value = 3.05
$ws.Range("L5").Value = 1.18
$ws.Range("N5").Value = 1.1
$ws.Range("O5").Value = 1.09
$ws.Range("P5").Value = 2.58
$ws.Range("Q5").Value = 1.09
$ws.Range("R5").Value = 1.76
$ws.Range("S5").Value = 1.76
$ws.Range("T6").Value = 1.42
$ws.Range("G8").Value = 2.5
$ws.Range("Q8").Value = 2.18
$ws.Range("S8").Value = 3.95
$ws.Range("F9").Value = 2.66
$ws.Range("H9").Value = 2.62
$ws.Range("I9").Value = 2.64
$ws.Range("P9").Value = 2.8
$ws.Range("L10").Value = 1.29
$ws.Range("S10").Value = 2.56
$ws.Range("T10").Value = 1.96
$ws.Range("AG10").Value = 34
$ws.Range("AH10").Value = 25
$ws.Range("F11").Value = 1.18
$ws.Range("G11").Value = 1.19
$ws.Range("H11").Value = 22
$ws.Range("I11").Value = 23
$ws.Range("J11").Value = 9
$ws.Range("K11").Value = 9.199999999999999
$ws.Range("N11").Value = 6.6
$ws.Range("P11").Value = 2.92
$ws.Range("Q11").Value = 1.49
$ws.Range("W11").Value = 6.2
$ws.Range("Y11").Value = 970
$ws.Range("AF11").Value = 7.8
$ws.Range("AG11").Value = 13
$ws.Range("AH11").Value = 140
$ws.Range("AN11").Value = 3.35
$ws.Range("K12").Value = 7.4
$ws.Range("P12").Value = 3.45
$ws.Range("S12").Value = 1.91
$ws.Range("AN12").Value = 3.4
$ws.Range("P13").Value = 2.14
$ws.Range("S13").Value = 3.1
$ws.Range("T13").Value = 1.9
$ws.Range("Q14").Value = 1.74
$ws.Range("T14").Value = 1.62
$ws.Range("AJ14").Value = 60
$ws.Range("F15").Value = 1.81
$ws.Range("G15").Value = 3.5
$ws.Range("H15").Value = 1.69
$ws.Range("I15").Value = 3.1
$ws.Range("J15").Value = 3.5
$ws.Range("S15").Value = 2.2
$ws.Range("T15").Value = 1.04
$ws.Range("U15").Value = 1.04
$ws.Range("V15").Value = 1.47
$ws.Range("W15").Value = 1.4
$ws.Range("F16").Value = 2.5
$ws.Range("Z16").Value = 22
$ws.Range("AI16").Value = 65
$ws.Range("AJ16").Value = 50
